# JORequestForm.xlsx - "Partial code for JOR"
#
# 1. Scroll the sheet so row 12 is the top visible row (was row 9).
# 2. J16 currently holds the formula =G16*I16 (cached 200) - replace it
#    with the plain static value 200 (drop the formula, keep the number).
# 3. K23:K26 currently hold hard-coded dates (45585-45588) - clear those
#    values while leaving the cell formatting (style) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JOR")

# --- 1. Scroll position: topLeftCell A9 -> A12 -------------------------
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1

# --- 2. J16: drop the formula, keep the resulting value ----------------
$ws.Range("J16").Value = 200

# --- 3. K23:K26: clear the hard-coded dates -----------------------------
$ws.Range("K23:K26").ClearContents()
